$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.915.80'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '1.642.49'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.53'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5094'
$ws.Range('E6').Value = '  +1.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.006'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2568'
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06394'
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07775'
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').Value = '1.653.42'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').Value = '0.0₅7844'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.74'
$ws.Range('E16').Value = '  +1.82%  '
$ws.Range('D17').Value = '25.975.42'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.006'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '197.82'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.437'
$ws.Range('E20').Value = '  +2.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.974'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.044'
$ws.Range('E22').Value = '  +1.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.007'
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('E24').Value = '  -2.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.81'
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1146'
$ws.Range('E26').Value = '  +1.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.898'
$ws.Range('E27').Value = '  +3.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.74'
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.239'
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05017'
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.190'
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.541'
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.364'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.8940'
$ws.Range('E35').Value = '  +0.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.592'
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('D37').Value = '1.133.94'
$ws.Range('E37').Value = '  -2.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5511'
$ws.Range('E38').Value = '  -0.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01561'
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.006'
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('B41').Value = 'mCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.554'
$ws.Range('E41').Value = '  -0.59%  '
$ws.Range('B42').Value = 'BabyDogeCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D42').Value = '0.0₈129'
$ws.Range('E42').Value = '  +11.07%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.638'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('E44').Value = '  +1.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.93'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('D46').Value = '1.780.15'
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4534'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.92'
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05089'
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.006'
$ws.Range('E51').Value = '  -0.04%  '
